$d = $word.ActiveDocument

# Title paragraph: merge "Questions:" " " "Completing" " " "the" " " "square"
# into a single run "Questions: Completing the square"
$d.Content.Find.Execute("Questions: Completing the square", $false, $false, $false, $false, $false, $true, 1, $false, "Questions: Completing the square", 2)

# Author paragraph: merge "Tom" " " "Coleman" into a single run "Tom Coleman"
$d.Content.Find.Execute("Tom Coleman", $false, $false, $false, $false, $false, $true, 1, $false, "Tom Coleman", 2)

# Abstract paragraph: merge the word-by-word runs into a single run
$d.Content.Find.Execute("A selection of questions for the study guide on completing the square.", $false, $false, $false, $false, $false, $true, 1, $false, "A selection of questions for the study guide on completing the square.", 2)
